$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update coin rows 2-51 (columns B,C,D,E) to match the latest scraped data.
# Values are set as explicit strings so Excel stores them as text, matching
# the inlineStr cell type used throughout this sheet (prices/links/percentages
# are display-formatted strings, not numeric cells).
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").Value = '96.325.12'
$ws.Range("E2").Value = '  +4.88%  '
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").Value = '3.620.82'
$ws.Range("E3").Value = '  +9.04%  '
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("B5").Value = 'Solana'
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D5").Value = '240.26'
$ws.Range("E5").Value = '  +4.39%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D6").Value = '639.59'
$ws.Range("E6").Value = '  +4.20%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value = '1.49'
$ws.Range("E7").Value = '  +6.34%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").Value = '0.401'
$ws.Range("E8").Value = '  +4.59%  '
$ws.Range("B9").Value = 'USDC'
$ws.Range("C9").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D9").Value = '0.999'
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("C10").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D10").Value = '1.01'
$ws.Range("E10").Value = '  +6.55%  '
$ws.Range("B11").Value = 'LidoStakedEther'
$ws.Range("C11").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("D11").Value = '3.617.63'
$ws.Range("E11").Value = '  +8.92%  '
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").Value = '43.33'
$ws.Range("E12").Value = '  +2.75%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").Value = '0.201'
$ws.Range("E13").Value = '  +3.78%  '
$ws.Range("B14").Value = 'Toncoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D14").Value = '6.35'
$ws.Range("E14").Value = '  +6.16%  '
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").Value = '4.312.63'
$ws.Range("E15").Value = '  +9.52%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").Value = '96.222.12'
$ws.Range("E16").Value = '  +5.05%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.0000254'
$ws.Range("E17").Value = '  +4.80%  '
$ws.Range("B18").Value = 'WrappedEther'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value = '3.625.55'
$ws.Range("E18").Value = '  +9.34%  '
$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").Value = '13.25'
$ws.Range("E19").Value = '  +22.24%  '
$ws.Range("B20").Value = 'Polkadot'
$ws.Range("C20").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D20").Value = '8.01'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("B21").Value = 'Chainlink'
$ws.Range("C21").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D21").Value = '18.14'
$ws.Range("E21").Value = '  +5.46%  '
$ws.Range("B22").Value = 'Stellar'
$ws.Range("C22").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D22").Value = '0.499'
$ws.Range("E22").Value = '  +11.37%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("C23").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D23").Value = '515.22'
$ws.Range("E23").Value = '  +4.92%  '
$ws.Range("B24").Value = 'SuiNetwork'
$ws.Range("C24").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D24").Value = '3.44'
$ws.Range("E24").Value = '  +0.41%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '0.0000197'
$ws.Range("E25").Value = '  +8.28%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = '6.70'
$ws.Range("E26").Value = '  +9.77%  '
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").Value = '97.42'
$ws.Range("E27").Value = '  +8.66%  '
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value = '3.816.83'
$ws.Range("E28").Value = '  +9.25%  '
$ws.Range("B29").Value = 'Aptos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").Value = '12.48'
$ws.Range("E29").Value = '  +5.98%  '
$ws.Range("B30").Value = 'PancakeSwap'
$ws.Range("C30").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D30").Value = '3.14'
$ws.Range("E30").Value = '  +21.15%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '11.56'
$ws.Range("E31").Value = '  +5.23%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").Value = '0.143'
$ws.Range("E32").Value = '  +3.40%  '
$ws.Range("B33").Value = 'Dai'
$ws.Range("C33").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D33").Value = '0.999'
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("B34").Value = 'Cronos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D34").Value = '0.181'
$ws.Range("E34").Value = '  +6.29%  '
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").Value = '0.990'
$ws.Range("E35").Value = '  -1.55%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '30.58'
$ws.Range("E36").Value = '  +8.96%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value = '0.568'
$ws.Range("E37").Value = '  +8.22%  '
$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").Value = '575.99'
$ws.Range("E38").Value = '  +4.04%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").Value = '7.86'
$ws.Range("E39").Value = '  +7.63%  '
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").Value = '1.47'
$ws.Range("E40").Value = '  +8.52%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("C41").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D41").Value = '0.152'
$ws.Range("E41").Value = '  +2.94%  '
$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.12%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("C43").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value = '0.925'
$ws.Range("E43").Value = '  +7.24%  '
$ws.Range("B44").Value = 'ImmutableX'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D44").Value = '1.73'
$ws.Range("E44").Value = '  +4.41%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = '0.0429'
$ws.Range("E45").Value = '  +5.39%  '
$ws.Range("B46").Value = 'WhiteBITCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D46").Value = '23.79'
$ws.Range("E46").Value = '  +0.47%  '
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Value = '5.70'
$ws.Range("E47").Value = '  +5.83%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '2.21'
$ws.Range("E48").Value = '  +5.64%  '
$ws.Range("B49").Value = 'MantraDAO'
$ws.Range("C49").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D49").Value = '3.52'
$ws.Range("E49").Value = '  -2.27%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").Value = '54.07'
$ws.Range("E50").Value = '  +4.51%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '8.15'
$ws.Range("E51").Value = '  +2.88%  '
